$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "Label" header in H1 (match style of existing header cells)
$ws.Range("H1").Value = "Label"
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

# Label values: 0 for Control rows, 1 for MDD rows (pattern repeats for rows 2-11 and 12-21)
$labels = @(0, 0, 0, 0, 0, 1, 1, 1, 1, 1)

for ($i = 0; $i -lt 10; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 8).Value = $labels[$i]
}

for ($i = 0; $i -lt 10; $i++) {
    $row = 12 + $i
    $ws.Cells.Item($row, 8).Value = $labels[$i]
}
